$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (GeoJSONStyleTestCase) comment: trim down the long text
$ws.Range("C25").Value = "Text looks wrong.  Too blocky."

# Row 10 (ScreenLabelsTestCase) comment: was "Multi-line not working"
$ws.Range("C10").Value = "Layout manager not quite right.  Could be layout size."

# Update the active selection to C14
$ws.Range("C14").Select()
